# Apply reordering of rows 2-12 (columns B & C) and clear column D (is_prefered)
# for the active worksheet, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order (row -> [id, speaker_variant]) for rows 2..12
$newValues = @(
    @("#de-hadriatische-zee", "De Hadriatische Zee"),
    @("#europe", "Europe"),
    @("#de-lydende-kerk", "De Lydende Kerk"),
    @("#kristekerk", "Kristekerk"),
    @("#gerechtigheid", "Gerechtigheid"),
    @("#de-rouw", "De Rouw"),
    @("#waarheid", "Waarheid"),
    @("#oostenryk", "Oostenryk"),
    @("#turkyen", "Turkyen"),
    @("#dapperheid", "Dapperheid"),
    @("#geweld", "Geweld")
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i][0]
    $ws.Cells.Item($row, 3).Value = $newValues[$i][1]
    $ws.Cells.Item($row, 4).Value = ""
}
